$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.598703456155553
$ws.Range("D2").Value = 0.1241508837947745

# Row 3
$ws.Range("C3").Value = 0.8888672110933737
$ws.Range("D3").Value = 0.3836906270191407

# Row 4
$ws.Range("C4").Value = -0.7815149157178162
$ws.Range("D4").Value = 0.4428294611047079

# Row 5
$ws.Range("C5").Value = 2.854872987114022
$ws.Range("D5").Value = 0.009209535467352037

# Row 6
$ws.Range("C6").Value = -0.5409637799467033
$ws.Range("D6").Value = 0.5939683560014086

# Row 7
$ws.Range("C7").Value = -1.460223059283481
$ws.Range("D7").Value = 0.1583630622201206

# Row 8
$ws.Range("C8").Value = 1.170027906105556
$ws.Range("D8").Value = 0.2545093056516181

# Row 9
$ws.Range("C9").Value = -1.568516489262696
$ws.Range("D9").Value = 0.1310327631987136

# Row 10
$ws.Range("C10").Value = 1.946681580761221
$ws.Range("D10").Value = 0.06445323060372909
$ws.Range("G10").Value = "No"

# Row 11
$ws.Range("C11").Value = 2.656418361799366
$ws.Range("D11").Value = 0.01441846804667901
